{"js": "// 1) Title: \"Facility Book System\" -> \"Facility Book System (Club)\"\n//    Added as a *separate* run (\" (Club)\") right after the existing run,\n//    so we use insertOoxml to keep the two runs distinct (a plain\n//    insertText would merge into the existing identically-formatted run).\nconst titleResults = context.document.body.search(\"Facility Book System\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  const titleRange = titleResults.items[0];\n  const clubOoxml =\n    '<?xml version=\"1.0\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/part.xml\" pkg:contentType=\"application/xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    '<w:r><w:t>Facility Book System</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> (Club)</w:t></w:r>' +\n    '</w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n  titleRange.insertOoxml(clubOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"After running the system, you can try the following features\"\n//    -> three runs: \"When\" + \" running the system, you can try the\n//    following features\" + \".\"\nconst featuresResults = context.document.body.search(\n  \"After running the system, you can try the following features\",\n  { matchCase: true }\n);\nfeaturesResults.load(\"items\");\nawait context.sync();\n\nif (featuresResults.items.length > 0) {\n  const featuresRange = featuresResults.items[0];\n  const featuresOoxml =\n    '<?xml version=\"1.0\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/part.xml\" pkg:contentType=\"application/xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    '<w:r><w:t>When</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> running the system, you can try the following features</w:t></w:r>' +\n    '<w:r><w:t>.</w:t></w:r>' +\n    '</w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n  featuresRange.insertOoxml(featuresOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Remove the stray \"_GoBack\" bookmark left after \"List and Cancel Bookings\".\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop script implementing the \"Club\" feature-doc refresh:\n#   1. Title \"Facility Book System\" gains a \" (Club)\" suffix (as its own run).\n#   2. \"After running the system, you can try the following features\" is\n#      rewritten as \"When running the system, you can try the following\n#      features.\" split across three runs (\"When\" / \" running the system,\n#      you can try the following features\" / \".\").\n#   3. The stray \"_GoBack\" bookmark left after \"List and Cancel Bookings\"\n#      is removed.\n\n$d = $word.ActiveDocument\n\n# -- 1) Title: add \" (Club)\" as a separate run right after the title text.\n$titleRange = $d.Content\n$found = $titleRange.Find.Execute(\"Facility Book System\")\nif ($found) {\n    # Clear the matched text first, then InsertXML at that (now-collapsed)\n    # position -- this engine's Range.InsertXML inserts alongside the range\n    # rather than replacing its contents, so clearing first keeps the\n    # result as two runs inside the original paragraph instead of adding a\n    # new one.\n    $titleRange.Text = \"\"\n    $titleXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/part.xml\" pkg:contentType=\"application/xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>Facility Book System</w:t></w:r><w:r><w:t xml:space=\"preserve\"> (Club)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $titleRange.InsertXML($titleXml)\n}\n\n# -- 2) \"After running...\" -> \"When running..., ... features.\"\n$featuresRange = $d.Content\n$found2 = $featuresRange.Find.Execute(\"After running the system, you can try the following features\")\nif ($found2) {\n    $featuresRange.Text = \"\"\n    $featuresXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/part.xml\" pkg:contentType=\"application/xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>When</w:t></w:r><w:r><w:t xml:space=\"preserve\"> running the system, you can try the following features</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $featuresRange.InsertXML($featuresXml)\n}\n\n# -- 3) Remove the leftover \"_GoBack\" bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
